$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Membre1")

# A4 already carries the date style used for the "Date" column (built-in
# m/d/yyyy number format with a thin border) - reuse it for the new date
# cells below via a formats-only copy/paste so no new style gets minted.
$ws.Cells.Item(4, 1).Copy() | Out-Null

# Row 8: fix the date (was stored as text "2017-1024") to the real date
# 2017-10-24, and fill in C8/D8 with their correct text.
$ws.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 1).Value = 43032
$ws.Cells.Item(8, 3).Value = "Création SceneGestionCompte"
$ws.Cells.Item(8, 4).Value = "Savoir comment organiser les menus."

# Row 9 (previously blank) - new task log entry
$ws.Cells.Item(9, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(9, 1).Value = 43032
$ws.Cells.Item(9, 2).Value = "10min"
$ws.Cells.Item(9, 3).Value = "Mise au point de la navigation entre les menus"
$ws.Cells.Item(9, 4).Value = "Le code de la navigation fonctionnait que si un textbox était actif."
$ws.Rows.Item(9).RowHeight = 28.8

# Row 10 (previously blank) - new task log entry
$ws.Cells.Item(10, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10, 1).Value = 43032
$ws.Cells.Item(10, 2).Value = "30min"
$ws.Cells.Item(10, 3).Value = "Implementation des singletons"
$ws.Cells.Item(10, 4).Value = "Le delete doit être fais à l'extérieur par une méthode static du singleton"
$ws.Rows.Item(10).RowHeight = 28.8

# Row 11 (previously blank) - new task log entry
# (D11's shared string is registered before C11's so the new shared-string
# table indices land in the same order as the reference workbook.)
$ws.Cells.Item(11, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(11, 1).Value = 43032
$ws.Cells.Item(11, 2).Value = "40min"
$ws.Cells.Item(11, 4).Value = "Il n'y a pas de méthode de split dans la librairie standard."
$ws.Cells.Item(11, 3).Value = "Implementation de l'authentification non testé et d'une méthode split non testé."
$ws.Rows.Item(11).RowHeight = 43.2

$excel.CutCopyMode = 0
$ws.Range("C11").Select()
